# [Add] XlsxToJson 변환 데이터 추가
#
# The source workbook was re-saved from Google Sheets format into a real
# Excel workbook and the data worksheet was renamed from the generic
# default "Sheet1" to the descriptive name "PlayerData" (matching the
# workbook's file name, PlayerData.xlsx). The underlying header/row data
# (player_hp, player_dam, player_mov, player_jump, smash_cooldown, and the
# single data row) is unchanged - only the sheet's display name changes.

$wb = $excel.ActiveWorkbook

# There is a single worksheet in this workbook; grab it defensively by
# name if present, otherwise fall back to the active sheet.
$ws = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Sheet1") {
        $ws = $sheet
        break
    }
}
if ($ws -eq $null) {
    $ws = $wb.ActiveSheet
}

$ws.Name = "PlayerData"
